$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 changes from text "-" to a numeric 0
$ws.Range("B2").Value = 0

# New header row entries
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# New "Area" column formulas (segment cross-sectional area), row 2 and 3 are
# standalone formulas, rows 4:15 form one shared/fill formula series
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Totals
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Update selection to match the authored state
$ws.Range("J2:K2").Select()

$excel.Calculate()
